# "List of parts updated" -- fill in new parts (rows 20-26) in the Main table,
# replacing the old placeholder rows (ToF module / Electric cubes / Tubular
# terminals) with the final, priced & sourced line items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use row 19 (a fully-formatted, already-priced row) as the formatting
# template for the new rows: A/B/C/H plain, D/E currency, F date, G hyperlink.
$ws.Range("A19:H19").Copy() | Out-Null
$ws.Range("A20:H20").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:H22").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:H23").PasteSpecial(-4122) | Out-Null
$ws.Range("A24:H24").PasteSpecial(-4122) | Out-Null
$ws.Range("A25:H25").PasteSpecial(-4122) | Out-Null
$ws.Range("A26:H26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$priceDate = 45173  # 2023-09-04

# Row 20: VL53L1X ToF module (ordered from AliExpress)
$ws.Range("A20").Value2 = "VL53L1X ToF module"
$ws.Range("B20").Value2 = 4
$ws.Range("C20").Value2 = "Ordered"
$ws.Range("D20").Value2 = 19.39
$ws.Range("F20").Value2 = $priceDate
$ws.Range("G20").Value2 = "https://pl.aliexpress.com/item/4000074204979.html"
$ws.Range("H20").Value2 = "AliExpress (SAMIORE Store)"
$ws.Hyperlinks.Add($ws.Range("G20"), "https://pl.aliexpress.com/item/4000074204979.html") | Out-Null
$ws.Range("G19").Copy() | Out-Null
$ws.Range("G20").PasteSpecial(-4122) | Out-Null

# Row 21: Security ON/OFF switch -- still order pending, keep its own
# (green-bordered) A21 style untouched, just refresh the text + date.
$ws.Range("F19").Copy() | Out-Null
$ws.Range("F21").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").Value2 = "Security ON/OFF switch"
$ws.Range("C21").Value2 = "Order pending"
$ws.Range("F21").Value2 = $priceDate

# Row 22: Electric cubes (12 x 2,5 mm2)
$ws.Range("A22").Value2 = "Electric cubes (12 x 2,5 mm2)"
$ws.Range("B22").Value2 = 3
$ws.Range("C22").Value2 = "Ordered"
$ws.Range("D22").Value2 = 1.57
$ws.Range("F22").Value2 = $priceDate
$ws.Range("G22").Value2 = "https://allegro.pl/oferta/listwa-zaciskowa-zlaczka-kostka-12-torowa-2-5mm2-8501387810"
$ws.Range("H22").Value2 = "Allegro (electro_24)"
$ws.Hyperlinks.Add($ws.Range("G22"), "https://allegro.pl/oferta/listwa-zaciskowa-zlaczka-kostka-12-torowa-2-5mm2-8501387810") | Out-Null
$ws.Range("G19").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null

# Row 23: Electric cubes (12 x 4 mm2) -- link pasted as plain text (no hyperlink)
$ws.Range("A23").Value2 = "Electric cubes (12 x 4 mm2)"
$ws.Range("B23").Value2 = 1
$ws.Range("C23").Value2 = "Ordered"
$ws.Range("D23").Value2 = 1.89
$ws.Range("F23").Value2 = $priceDate
$ws.Range("G23").Value2 = "https://allegro.pl/oferta/zlaczka-listwa-instalacyjna-kablowa-12-torowa-4mm-7270241593"
$ws.Range("G23").ClearFormats()
$ws.Range("H23").Value2 = "Allegro (electro_24)"

# Row 24: Tubular terminals set -- link pasted as plain text (no hyperlink)
$ws.Range("A24").Value2 = "Tubular terminals set"
$ws.Range("B24").Value2 = 1
$ws.Range("C24").Value2 = "Ordered"
$ws.Range("D24").Value2 = 54.9
$ws.Range("F24").Value2 = $priceDate
$ws.Range("G24").Value2 = "https://allegro.pl/oferta/zaciskarka-do-konektorow-koncowek-500-tulejek-13993749597"
$ws.Range("G24").ClearFormats()
$ws.Range("H24").Value2 = "Allegro (electro_24)"

# Row 25: Universal PCB set
$ws.Range("A25").Value2 = "Universal PCB set"
$ws.Range("B25").Value2 = 1
$ws.Range("C25").Value2 = "Ordered"
$ws.Range("D25").Value2 = 33.99
$ws.Range("F25").Value2 = $priceDate
$ws.Range("G25").Value2 = "https://www.amazon.pl/gp/product/B07V25W5RT"
$ws.Range("H25").Value2 = "Amazon (AZDelivery)"
$ws.Hyperlinks.Add($ws.Range("G25"), "https://www.amazon.pl/gp/product/B07V25W5RT") | Out-Null
$ws.Range("G19").Copy() | Out-Null
$ws.Range("G25").PasteSpecial(-4122) | Out-Null

# Row 26: IR module set
$ws.Range("A26").Value2 = "IR module set"
$ws.Range("B26").Value2 = 1
$ws.Range("C26").Value2 = "Ordered"
$ws.Range("D26").Value2 = 20.49
$ws.Range("F26").Value2 = $priceDate
$ws.Range("G26").Value2 = "https://www.amazon.pl/gp/product/B07V9XD2R6"
$ws.Range("H26").Value2 = "Amazon (AZDelivery)"
$ws.Hyperlinks.Add($ws.Range("G26"), "https://www.amazon.pl/gp/product/B07V9XD2R6") | Out-Null
$ws.Range("G19").Copy() | Out-Null
$ws.Range("G26").PasteSpecial(-4122) | Out-Null

# Match the author's last active selection.
$ws.Range("I31").Select() | Out-Null
